# Generate Report for Handback
#
# Populates the "Latest Target File" / "Latest Handback File" / "Latest
# Handback DateTime" columns on the zh-cn and de-de report sheets now that
# both files have been handed back and are in sync with en-US, and updates
# the Overview sheet's status text + a handful of column widths so the new,
# longer values are readable.

$wb = $excel.ActiveWorkbook

$baseUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a61c8a44fc75e12e62856c77574347614dc42838/e2e/"

$adName    = "ad5cd1d7-6aed-4e16-9cfb-44a5f403785b.md"
$d1Name    = "d12c01b6-d014-4d19-bc48-bb773c602721.md"

$adZhXlf   = "ad5cd1d7-6aed-4e16-9cfb-44a5f403785b.ea78ba7872ca846d39e6917886e44267088ffac4.zh-cn.xlf"
$d1ZhXlf   = "d12c01b6-d014-4d19-bc48-bb773c602721.fd690c00a51446a93947cd4ba240d7c4e52a69f5.zh-cn.xlf"
$adDeXlf   = "ad5cd1d7-6aed-4e16-9cfb-44a5f403785b.ea78ba7872ca846d39e6917886e44267088ffac4.de-de.xlf"
$d1DeXlf   = "d12c01b6-d014-4d19-bc48-bb773c602721.fd690c00a51446a93947cd4ba240d7c4e52a69f5.de-de.xlf"

$handbackDateTimeZhCn = "2017-02-09 17:38:18"
$handbackDateTimeDeDe = "2017-02-09 17:38:46"

# ---------------------------------------------------------------------
# Overview sheet: status text now reflects the handback, and the status
# columns need to be wider to fit it.
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"
$overview.Columns.Item(5).ColumnWidth = 29.9777050018311
$overview.Columns.Item(6).ColumnWidth = 29.9777050018311

# ---------------------------------------------------------------------
# zh-cn report sheet
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = 29.9777050018311
$zhcn.Columns.Item(10).ColumnWidth = 40
$zhcn.Columns.Item(11).ColumnWidth = 40

$zhcn.Range("C2").Value = "Handed back: in sync with en-US"
$zhcn.Range("C3").Value = "Handed back: in sync with en-US"

$zhcn.Range("K2").Value = $adZhXlf
$zhcn.Range("L2").Value = $handbackDateTimeZhCn
$zhcn.Hyperlinks.Add($zhcn.Range("J2"), $baseUrl + $adName, $null, $null, $adName) | Out-Null

$zhcn.Range("K3").Value = $d1ZhXlf
$zhcn.Range("L3").Value = $handbackDateTimeZhCn
$zhcn.Hyperlinks.Add($zhcn.Range("J3"), $baseUrl + $d1Name, $null, $null, $d1Name) | Out-Null

# ---------------------------------------------------------------------
# de-de report sheet
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = 29.9777050018311
$dede.Columns.Item(10).ColumnWidth = 40
$dede.Columns.Item(11).ColumnWidth = 40

$dede.Range("C2").Value = "Handed back: in sync with en-US"
$dede.Range("C3").Value = "Handed back: in sync with en-US"

$dede.Range("K2").Value = $adDeXlf
$dede.Range("L2").Value = $handbackDateTimeDeDe
$dede.Hyperlinks.Add($dede.Range("J2"), $baseUrl + $adName, $null, $null, $adName) | Out-Null

$dede.Range("K3").Value = $d1DeXlf
$dede.Range("L3").Value = $handbackDateTimeDeDe
$dede.Hyperlinks.Add($dede.Range("J3"), $baseUrl + $d1Name, $null, $null, $d1Name) | Out-Null
